# Short course MDR-TB regimen improves outcomes
# Insert a new parameter row for "program_prop_treatment_success_shortcoursemdr"
# into the "constants" sheet, just above the existing
# "program_timeperiod_acf_rounds" row (old row 48), pushing all rows
# below it down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new blank row at row 48 (shifts rows 48-118 down to 49-119,
# inheriting formatting/styles from the row above, as Excel normally does).
$ws.Rows.Item(48).Insert()

# Populate the new row with the new parameter data.
$ws.Cells.Item(48, 1).Value = "program_prop_treatment_success_shortcoursemdr"
$ws.Cells.Item(48, 2).Value = 0.879
$ws.Cells.Item(48, 5).Value = "Treatment success under shortcourse MDR-TB regimens"

# Update the sheet's selection to reflect where the author left off editing.
$null = $ws.Activate()
$ws.Range("B49").Select() | Out-Null
